$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the existing row 100 (old rows 100:120
# shift down to become rows 103:123, matching the new dimension A1:R123).
$ws.Rows("100:102").Insert()

# --- New row 100 (Calidad: Extra) ---
$ws.Range("A100").Value2 = 7
$ws.Range("B100").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value2 = "Ñuble"
$ws.Range("D100").Value2 = 44559
$ws.Range("E100").Value2 = 16
$ws.Range("F100").Value2 = 100112028
$ws.Range("G100").Value2 = "Sandia"
$ws.Range("H100").Value2 = "Sin especificar"
$ws.Range("I100").Value2 = "Extra"
$ws.Range("J100").Value2 = 300
$ws.Range("K100").Value2 = 2500
$ws.Range("L100").Value2 = 2500
$ws.Range("M100").Value2 = 2500
$ws.Range("N100").Value2 = "$/unidad"
$ws.Range("O100").Value2 = "Región de O'Higgins"
$ws.Range("P100").Value2 = 2500
$ws.Range("Q100").Value2 = 1
$ws.Range("R100").Value2 = "Hortaliza"

# --- New row 101 (Calidad: Primera) ---
$ws.Range("A101").Value2 = 7
$ws.Range("B101").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C101").Value2 = "Ñuble"
$ws.Range("D101").Value2 = 44559
$ws.Range("E101").Value2 = 16
$ws.Range("F101").Value2 = 100112028
$ws.Range("G101").Value2 = "Sandia"
$ws.Range("H101").Value2 = "Sin especificar"
$ws.Range("I101").Value2 = "Primera"
$ws.Range("J101").Value2 = 400
$ws.Range("K101").Value2 = 2000
$ws.Range("L101").Value2 = 2200
$ws.Range("M101").Value2 = 2100
$ws.Range("N101").Value2 = "$/unidad"
$ws.Range("O101").Value2 = "Región de O'Higgins"
$ws.Range("P101").Value2 = 2100
$ws.Range("Q101").Value2 = 1
$ws.Range("R101").Value2 = "Hortaliza"

# --- New row 102 (Calidad: Segunda) ---
$ws.Range("A102").Value2 = 7
$ws.Range("B102").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C102").Value2 = "Ñuble"
$ws.Range("D102").Value2 = 44559
$ws.Range("E102").Value2 = 16
$ws.Range("F102").Value2 = 100112028
$ws.Range("G102").Value2 = "Sandia"
$ws.Range("H102").Value2 = "Sin especificar"
$ws.Range("I102").Value2 = "Segunda"
$ws.Range("J102").Value2 = 300
$ws.Range("K102").Value2 = 1800
$ws.Range("L102").Value2 = 1800
$ws.Range("M102").Value2 = 1800
$ws.Range("N102").Value2 = "$/unidad"
$ws.Range("O102").Value2 = "Región de O'Higgins"
$ws.Range("P102").Value2 = 1800
$ws.Range("Q102").Value2 = 1
$ws.Range("R102").Value2 = "Hortaliza"
